$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1340.18
$ws.Range("J17").Value = 1354.409
$ws.Range("L17").Value = 4063.227
$ws.Range("N17").Value = -4399.227000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 64950.625
$ws.Range("J51").Value = 127366
$ws.Range("L51").Value = 127366
$ws.Range("N51").Value = -128334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2930.4666
$ws.Range("I53").Value = 3827.9092
$ws.Range("J53").Value = 462.5
$ws.Range("K53").Value = 3827.9092
$ws.Range("L53").Value = 462.5
$ws.Range("M53").Value = -3190.9092
$ws.Range("N53").Value = -1736.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 181800.47
$ws.Range("I138").Value = 455521.88
$ws.Range("J138").Value = 5402.2446
$ws.Range("K138").Value = 1366565.64
$ws.Range("L138").Value = 16206.7338
$ws.Range("M138").Value = -1361425.64
$ws.Range("N138").Value = -26486.7338

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("M19").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6205.6
$ws.Range("I32").Value = 5511.1577
$ws.Range("J32").Value = 19400
$ws.Range("K32").Value = 5511.1577
$ws.Range("L32").Value = 19400
$ws.Range("M32").Value = -5224.1577
$ws.Range("N32").Value = -19974

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4052.849
$ws.Range("I74").Value = 40666.668
$ws.Range("J74").Value = 1856.02
$ws.Range("K74").Value = 40666.668
$ws.Range("L74").Value = 1856.02
$ws.Range("M74").Value = -39792.668
$ws.Range("N74").Value = -3604.02

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4052.849
$ws.Range("I77").Value = 40666.668
$ws.Range("J77").Value = 1856.02
$ws.Range("K77").Value = 203333.34
$ws.Range("L77").Value = 9280.1
$ws.Range("M77").Value = -198965.34
$ws.Range("N77").Value = -18016.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 89999.5
$ws.Range("I109").Value = 90000
$ws.Range("J109").Value = 89999
$ws.Range("K109").Value = 90000
$ws.Range("L109").Value = 89999
$ws.Range("M109").Value = -88613
$ws.Range("N109").Value = -92773

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8041.857
$ws.Range("J86").Value = 3910
$ws.Range("L86").Value = 3910
$ws.Range("N86").Value = -6156

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 8041.857
$ws.Range("J89").Value = 3910
$ws.Range("L89").Value = 19550
$ws.Range("N89").Value = -30782

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 82399.8
$ws.Range("J132").Value = 82399.8
$ws.Range("L132").Value = 82399.8
$ws.Range("N132").Value = -92519.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 215.75
$ws.Range("J22").Value = 220.66667
$ws.Range("L22").Value = 220.66667
$ws.Range("N22").Value = -920.6666700000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4759.3335
$ws.Range("I31").Value = 1294.5
$ws.Range("J31").Value = 5749.2856
$ws.Range("K31").Value = 1294.5
$ws.Range("L31").Value = 5749.2856
$ws.Range("M31").Value = -999.5
$ws.Range("N31").Value = -6339.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4759.3335
$ws.Range("I34").Value = 1294.5
$ws.Range("J34").Value = 5749.2856
$ws.Range("K34").Value = 1294.5
$ws.Range("L34").Value = 5749.2856
$ws.Range("M34").Value = -1092.5
$ws.Range("N34").Value = -6153.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9881.4
$ws.Range("I62").Value = 6660.6
$ws.Range("J62").Value = 13102.2
$ws.Range("K62").Value = 6660.6
$ws.Range("L62").Value = 13102.2
$ws.Range("M62").Value = -6036.6
$ws.Range("N62").Value = -14350.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 9881.4
$ws.Range("I65").Value = 6660.6
$ws.Range("J65").Value = 13102.2
$ws.Range("K65").Value = 33303
$ws.Range("L65").Value = 65511
$ws.Range("M65").Value = -30183
$ws.Range("N65").Value = -71751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5534382
$ws.Range("I99").Value = 8935925
$ws.Range("J99").Value = 6875
$ws.Range("K99").Value = 8935925
$ws.Range("L99").Value = 6875
$ws.Range("M99").Value = -8934427
$ws.Range("N99").Value = -9871

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5534382
$ws.Range("I126").Value = 8935925
$ws.Range("J126").Value = 6875
$ws.Range("K126").Value = 26807775
$ws.Range("L126").Value = 20625
$ws.Range("M126").Value = -26805305
$ws.Range("N126").Value = -25565

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 409556.62
$ws.Range("J141").Value = 515094
$ws.Range("L141").Value = 515094
$ws.Range("N141").Value = -525454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 125.1579
$ws.Range("I7").Value = 120.5625
$ws.Range("K7").Value = 361.6875
$ws.Range("M7").Value = -249.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1631.7142
$ws.Range("J60").Value = 3533.3333
$ws.Range("L60").Value = 10599.9999
$ws.Range("N60").Value = -11101.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 100166000
$ws.Range("I80").Value = 250006000
$ws.Range("J80").Value = 272666.66
$ws.Range("K80").Value = 750018000
$ws.Range("L80").Value = 817999.98
$ws.Range("M80").Value = -750017064
$ws.Range("N80").Value = -819871.98

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 100166000
$ws.Range("I83").Value = 250006000
$ws.Range("J83").Value = 272666.66
$ws.Range("K83").Value = 2250054000
$ws.Range("L83").Value = 2453999.94
$ws.Range("M83").Value = -2250049320
$ws.Range("N83").Value = -2463359.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 4872.6
$ws.Range("I115").Value = 3994
$ws.Range("K115").Value = 11982
$ws.Range("M115").Value = -10807

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1824.5834
$ws.Range("I121").Value = 1038.8
$ws.Range("J121").Value = 2385.8572
$ws.Range("K121").Value = 3116.4
$ws.Range("L121").Value = 7157.571599999999
$ws.Range("M121").Value = -1806.4
$ws.Range("N121").Value = -9777.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4830.4585
$ws.Range("J131").Value = 1974.4615
$ws.Range("L131").Value = 5923.3845
$ws.Range("N131").Value = -16003.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 63325.3
$ws.Range("I132").Value = 3033.3333
$ws.Range("J132").Value = 89164.71000000001
$ws.Range("K132").Value = 27299.9997
$ws.Range("L132").Value = 802482.39
$ws.Range("M132").Value = -24769.9997
$ws.Range("N132").Value = -807542.39

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 40000
$ws.Range("I113").Value = 110000
$ws.Range("K113").Value = 110000
$ws.Range("M113").Value = -107830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18830.049
$ws.Range("I7").Value = 22743.732
$ws.Range("K7").Value = 22743.732
$ws.Range("M7").Value = -22631.732

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 25491.428
$ws.Range("I40").Value = 35495.46
$ws.Range("J40").Value = 16821.268
$ws.Range("K40").Value = 35495.46
$ws.Range("L40").Value = 16821.268
$ws.Range("M40").Value = -35359.46
$ws.Range("N40").Value = -17093.268

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3701.9473
$ws.Range("I93").Value = 4476.8076
$ws.Range("K93").Value = 4476.8076
$ws.Range("M93").Value = -3228.8076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 12800
$ws.Range("J100").Value = 9750
$ws.Range("L100").Value = 9750
$ws.Range("N100").Value = -10832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 18830.049
$ws.Range("I126").Value = 22743.732
$ws.Range("K126").Value = 68231.196
$ws.Range("M126").Value = -65761.196

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 100005750
$ws.Range("I2").Value = 125007040
$ws.Range("K2").Value = 125007040
$ws.Range("M2").Value = -125006928

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 28074.75
$ws.Range("J126").Value = 5399
$ws.Range("L126").Value = 16197
$ws.Range("N126").Value = -21137

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4191.84
$ws.Range("I132").Value = 4366.8
$ws.Range("J132").Value = 3710.7
$ws.Range("K132").Value = 13100.4
$ws.Range("L132").Value = 11132.1
$ws.Range("M132").Value = -10570.4
$ws.Range("N132").Value = -16192.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3852903.5
$ws.Range("I136").Value = 5134871.5
$ws.Range("K136").Value = 15404614.5
$ws.Range("M136").Value = -15402064.5
